$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the formula in A2: 20*90*24 -> 1*90*24
$ws.Range("A2").Formula = "=1*90*24"

# Recalculate so dependent formulas (F2:F7) pick up the new value
$excel.Calculate()

# Update the selected cell to A3 (as reflected in the saved sheet view)
$ws.Range("A3").Select()
